$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Resolving-Mac"
$ws.Cells.Item(2, 3).Value = "Ccl12"
$ws.Cells.Item(2, 4).Value = "Ccr5"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 20.23247666666667
$ws.Cells.Item(2, 8).Value = 60.69743
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.022105
$ws.Cells.Item(2, 14).Value = 0.066315
$ws.Cells.Item(2, 15).Value = 0.0007557226718989593
$ws.Cells.Item(2, 16).Value = 0.0007557226718989592
$ws.Cells.Item(2, 17).Value = 0.4472388967166667
$ws.Cells.Item(2, 18).Value = 4.02515007045
$ws.Cells.Item(2, 19).Value = 0.0007557226718989593
$ws.Cells.Item(2, 20).Value = 0.0007557226718989592

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Resolving-Mac"
$ws.Cells.Item(3, 3).Value = "Ccl12"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 20.23247666666667
$ws.Cells.Item(3, 8).Value = 60.69743
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.4010506666666667
$ws.Cells.Item(3, 14).Value = 1.203152
$ws.Cells.Item(3, 15).Value = 0.01371106452749117
$ws.Cells.Item(3, 16).Value = 0.01371106452749117
$ws.Cells.Item(3, 17).Value = 8.114248255484444
$ws.Cells.Item(3, 18).Value = 73.02823429936
$ws.Cells.Item(3, 19).Value = 0.01371106452749117
$ws.Cells.Item(3, 20).Value = 0.01371106452749117

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Resolving-Mac"
$ws.Cells.Item(4, 3).Value = "Ccl12"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 20.23247666666667
$ws.Cells.Item(4, 8).Value = 60.69743
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 28.82699233333333
$ws.Cells.Item(4, 14).Value = 86.480977
$ws.Cells.Item(4, 15).Value = 0.9855332128006099
$ws.Cells.Item(4, 16).Value = 0.9855332128006098
$ws.Cells.Item(4, 17).Value = 583.2414497543456
$ws.Cells.Item(4, 18).Value = 5249.17304778911
$ws.Cells.Item(4, 19).Value = 0.9855332128006099
$ws.Cells.Item(4, 20).Value = 0.9855332128006098

# Remove now-obsolete rows 5-7
$ws.Range("A5:T7").EntireRow.Delete()

Write-Output "done"